# ---------------------------------------------------------------------------
# Edit: add UOM(S) master data derived "Concrete Delivered Qty" column to the
# Dispatch Plant Day and Dispatch Plant Month sheets, refresh the "Generated"
# timestamp on the Cover sheet, and let the dependent Dispatch vs AR sheet's
# shared-string references follow automatically.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Cover sheet: bump the "Generated" timestamp.
# ---------------------------------------------------------------------------
$cover = $wb.Worksheets.Item("Cover")
$cover.Range("B4").Value = "2026-02-17 10:55"

# ---------------------------------------------------------------------------
# 2) Dispatch Plant Day: insert a new "Concrete Delivered Qty" column
#    (column E) between "Delivered Qty" (D) and "Revenue" (old E, now F).
# ---------------------------------------------------------------------------
$day = $wb.Worksheets.Item("Dispatch Plant Day")

$day.Columns.Item(5).Insert()
$day.Columns.Item(5).ColumnWidth = 21.25

$day.Range("E1").Value = "Concrete Delivered Qty"

$dayValues = @{
    2  = 0
    3  = 108.5
    4  = 71.5
    5  = 216.25
    6  = 0
    7  = 235
    8  = 124
    9  = 183.25
    10 = 436.75
    11 = 0
    12 = 0
    13 = 0
    14 = 409.25
    15 = 56.5
    16 = 9
    17 = 44
    18 = 39
    19 = 45.5
    20 = 102
    21 = 128
    22 = 0
    23 = 75
    24 = 105
    25 = 226
    26 = 362
    27 = 0
    28 = 344.75
    29 = 112.5
    30 = 433.25
    31 = 657.25
    32 = 0
    33 = 282
    34 = 54.5
    35 = 110.5
    36 = 126.5
    37 = 0
    38 = 159
    39 = 111.25
    40 = 51.5
    41 = 280.75
    42 = 0
    43 = 0
    44 = 155.5
    45 = 23
    46 = 46.5
    47 = 323.5
    48 = 0
    49 = 343
    50 = 109.5
    51 = 31.5
    52 = 496.25
    53 = 535.25
}

foreach ($row in 2..53) {
    $day.Range("E$row").Value = $dayValues[$row]
}

# ---------------------------------------------------------------------------
# 3) Dispatch Plant Month: insert the matching "Concrete Delivered Qty"
#    column (column F) between "Delivered Qty" (E) and "Revenue" (old F,
#    now G).
# ---------------------------------------------------------------------------
$month = $wb.Worksheets.Item("Dispatch Plant Month")

$month.Columns.Item(6).Insert()
$month.Columns.Item(6).ColumnWidth = 21.25

$month.Range("F1").Value = "Concrete Delivered Qty"

$monthValues = @{
    2 = 802.75
    3 = 1300
    4 = 3083.75
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 2577.75
}

foreach ($row in 2..9) {
    $month.Range("F$row").Value = $monthValues[$row]
}
